$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 25643196
$ws.Range("I100").Value = 2256.5557
$ws.Range("J100").Value = 83335310
$ws.Range("K100").Value = 2256.5557
$ws.Range("L100").Value = 83335310
$ws.Range("M100").Value = -1715.5557
$ws.Range("N100").Value = -83336392
$ws.Range("H125").Value = 2526.889
$ws.Range("I125").Value = 800.6667
$ws.Range("J125").Value = 3390
$ws.Range("K125").Value = 7206.0003
$ws.Range("L125").Value = 30510
$ws.Range("M125").Value = -4746.0003
$ws.Range("N125").Value = -35430
$ws.Range("H129").Value = 1502
$ws.Range("I129").Value = 483.33334
$ws.Range("J129").Value = 1640.909
$ws.Range("K129").Value = 1450.00002
$ws.Range("L129").Value = 4922.727000000001
$ws.Range("M129").Value = 3549.99998
$ws.Range("N129").Value = -14922.727
$ws.Range("H137").Value = 2675.1875
$ws.Range("I137").Value = 1550
$ws.Range("J137").Value = 3350.3
$ws.Range("K137").Value = 4650
$ws.Range("L137").Value = 10050.9
$ws.Range("M137").Value = -2100
$ws.Range("N137").Value = -15150.9
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H80").Value = 36055
$ws.Range("J80").Value = 36055
$ws.Range("L80").Value = 36055
$ws.Range("N80").Value = -38051
$ws.Range("H83").Value = 36055
$ws.Range("J83").Value = 36055
$ws.Range("L83").Value = 108165
$ws.Range("N83").Value = -118149
$ws.Range("H122").Value = 2201.5
$ws.Range("I122").Value = 2036.6428
$ws.Range("J122").Value = 2778.5
$ws.Range("K122").Value = 6109.928400000001
$ws.Range("L122").Value = 8335.5
$ws.Range("M122").Value = -3659.928400000001
$ws.Range("N122").Value = -13235.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 39982.668
$ws.Range("J35").Value = 39982.668
$ws.Range("L35").Value = 39982.668
$ws.Range("N35").Value = -40602.668
$ws.Range("H60").Value = 20050
$ws.Range("J60").Value = 20050
$ws.Range("L60").Value = 20050
$ws.Range("N60").Value = -21248
$ws.Range("H82").Value = 15174.714
$ws.Range("I82").Value = 5851.4
$ws.Range("J82").Value = 38483
$ws.Range("K82").Value = 5851.4
$ws.Range("L82").Value = 38483
$ws.Range("M82").Value = -5468.4
$ws.Range("N82").Value = -39249
$ws.Range("H85").Value = 15174.714
$ws.Range("I85").Value = 5851.4
$ws.Range("J85").Value = 38483
$ws.Range("K85").Value = 5851.4
$ws.Range("L85").Value = 38483
$ws.Range("M85").Value = -4525.4
$ws.Range("N85").Value = -41135
$ws.Range("H94").Value = 1393.1333
$ws.Range("I94").Value = 1569.5
$ws.Range("J94").Value = 1275.5555
$ws.Range("K94").Value = 1569.5
$ws.Range("L94").Value = 1275.5555
$ws.Range("M94").Value = -1118.5
$ws.Range("N94").Value = -2177.5555
$ws.Range("H125").Value = 42500
$ws.Range("J125").Value = 42500
$ws.Range("L125").Value = 42500
$ws.Range("N125").Value = -52340
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 37899
$ws.Range("J51").Value = 37899
$ws.Range("L51").Value = 37899
$ws.Range("N51").Value = -39371
$ws.Range("H58").Value = 1647.7391
$ws.Range("I58").Value = 1805.742
$ws.Range("K58").Value = 1805.742
$ws.Range("M58").Value = -1602.742
$ws.Range("H61").Value = 37899
$ws.Range("J61").Value = 37899
$ws.Range("L61").Value = 37899
$ws.Range("N61").Value = -38595
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H136").Value = 1647.7391
$ws.Range("I136").Value = 1805.742
$ws.Range("K136").Value = 5417.226
$ws.Range("M136").Value = -2867.226
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 825
$ws.Range("I5").Value = 825
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 2475
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -2363
$ws.Range("N5").ClearContents()
$ws.Range("H135").Value = 825
$ws.Range("I135").Value = 825
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 7425
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -4890
$ws.Range("N135").Value = -4890
$ws.Range("H137").Value = 1913.8462
$ws.Range("I137").Value = 1263.3334
$ws.Range("J137").Value = 2471.4285
$ws.Range("K137").Value = 3790.0002
$ws.Range("L137").Value = 7414.2855
$ws.Range("M137").Value = 1309.9998
$ws.Range("N137").Value = -17614.2855
$ws.Range("H138").Value = 40001612
$ws.Range("I138").Value = 76923944
$ws.Range("J138").Value = 2415
$ws.Range("K138").Value = 230771832
$ws.Range("L138").Value = 7245
$ws.Range("M138").Value = -230766692
$ws.Range("N138").Value = -17525
$ws.Range("H139").Value = 29741
$ws.Range("I139").Value = 1438.0741
$ws.Range("K139").Value = 4314.2223
$ws.Range("M139").Value = 825.7776999999996
$ws.Range("H141").Value = 13540.833
$ws.Range("I141").Value = 7898
$ws.Range("J141").Value = 17571.428
$ws.Range("K141").Value = 23694
$ws.Range("L141").Value = 52714.284
$ws.Range("M141").Value = -18514
$ws.Range("N141").Value = -63074.284
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H45").Value = 14612.5
$ws.Range("J45").Value = 14612.5
$ws.Range("L45").Value = 14612.5
$ws.Range("N45").Value = -15730.5
$ws.Range("H51").Value = 14763.333
$ws.Range("J51").Value = 14763.333
$ws.Range("L51").Value = 14763.333
$ws.Range("N51").Value = -15781.333
$ws.Range("H80").Value = 2516.6956
$ws.Range("I80").Value = 2872.4
$ws.Range("J80").Value = 2243.077
$ws.Range("K80").Value = 2872.4
$ws.Range("L80").Value = 2243.077
$ws.Range("M80").Value = -1874.4
$ws.Range("N80").Value = -4239.077
$ws.Range("H83").Value = 2516.6956
$ws.Range("I83").Value = 2872.4
$ws.Range("J83").Value = 2243.077
$ws.Range("K83").Value = 14362
$ws.Range("L83").Value = 11215.385
$ws.Range("M83").Value = -9370
$ws.Range("N83").Value = -21199.385
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
$ws.Range("H102").Value = 1528.9375
$ws.Range("I102").Value = 1230.8667
$ws.Range("K102").Value = 1230.8667
$ws.Range("M102").Value = 391.1333
$ws.Range("H122").Value = 3474.353
$ws.Range("I122").Value = 3183.1428
$ws.Range("J122").Value = 4833.3335
$ws.Range("K122").Value = 9549.428400000001
$ws.Range("L122").Value = 14500.0005
$ws.Range("M122").Value = -7099.428400000001
$ws.Range("N122").Value = -19400.0005
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3083.9524
$ws.Range("I122").Value = 2320.3333
$ws.Range("J122").Value = 3656.6667
$ws.Range("K122").Value = 6960.999899999999
$ws.Range("L122").Value = 10970.0001
$ws.Range("M122").Value = -4510.999899999999
$ws.Range("N122").Value = -15870.0001
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7484.5
$ws.Range("I62").Value = 5002
$ws.Range("J62").Value = 7981
$ws.Range("K62").Value = 5002
$ws.Range("L62").Value = 7981
$ws.Range("M62").Value = -4378
$ws.Range("N62").Value = -9229
$ws.Range("H65").Value = 7484.5
$ws.Range("I65").Value = 5002
$ws.Range("J65").Value = 7981
$ws.Range("K65").Value = 25010
$ws.Range("L65").Value = 39905
$ws.Range("M65").Value = -21890
$ws.Range("N65").Value = -46145
$ws.Range("H96").Value = 2005.3334
$ws.Range("I96").Value = 1564
$ws.Range("K96").Value = 1564
$ws.Range("M96").Value = -191
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()
$ws.Range("H106").Value = 39886.5
$ws.Range("J106").Value = 39886.5
$ws.Range("L106").Value = 39886.5
$ws.Range("N106").Value = -42410.5
$ws.Range("H109").Value = 29999
$ws.Range("J109").Value = 29999
$ws.Range("L109").Value = 29999
$ws.Range("N109").Value = -32773
$ws.Range("H113").Value = 39429.117
$ws.Range("I113").Value = 55737.445
$ws.Range("K113").Value = 167212.335
$ws.Range("M113").Value = -165042.335
$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()
$ws.Range("H118").Value = 54916
$ws.Range("J118").Value = 54916
$ws.Range("L118").Value = 54916
$ws.Range("N118").Value = -58230
$ws.Range("H122").Value = 2735.25
$ws.Range("I122").Value = 3057.1428
$ws.Range("J122").Value = 1984.1666
$ws.Range("K122").Value = 9171.428400000001
$ws.Range("L122").Value = 5952.4998
$ws.Range("M122").Value = -6721.428400000001
$ws.Range("N122").Value = -10852.4998
$ws.Range("H133").Value = 63985
$ws.Range("J133").Value = 63985
$ws.Range("L133").Value = 63985
$ws.Range("N133").Value = -74105
